# Re-apply the latest crypto snapshot values scraped by the GitHub Action.
# Cell D/E values are plain text (prices/volume strings), not numbers, in the
# source workbook (t="inlineStr"). A leading apostrophe forces Excel to keep
# them as literal text instead of silently re-parsing them as numbers/dates
# (which would strip things like trailing zeros or multi-dot price strings).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'28.477.52"
$ws.Range("E2").Value = "'  -0.93%  "
$ws.Range("D3").Value = "'1.878.05"
$ws.Range("E3").Value = "'  -1.71%  "
$ws.Range("D4").Value = "'1.014"
$ws.Range("E4").Value = "'  -1.40%  "
$ws.Range("D5").Value = "'316.42"
$ws.Range("E5").Value = "'  -1.38%  "
$ws.Range("D6").Value = "'1.012"
$ws.Range("E6").Value = "'  -1.74%  "
$ws.Range("D7").Value = "'0.5110"
$ws.Range("E7").Value = "'  -2.06%  "
$ws.Range("D8").Value = "'0.3959"
$ws.Range("E8").Value = "'  -0.06%  "
$ws.Range("D9").Value = "'0.08409"
$ws.Range("E9").Value = "'  +0.17%  "
$ws.Range("D10").Value = "'1.110"
$ws.Range("E10").Value = "'  -2.41%  "
$ws.Range("D11").Value = "'6.271"
$ws.Range("E11").Value = "'  -0.85%  "
$ws.Range("D12").Value = "'1.887.34"
$ws.Range("E12").Value = "'  -1.90%  "
$ws.Range("D13").Value = "'20.48"
$ws.Range("E13").Value = "'  -1.12%  "
$ws.Range("D14").Value = "'7.276"
$ws.Range("E14").Value = "'  -0.92%  "
$ws.Range("D15").Value = "'1.014"
$ws.Range("E15").Value = "'  -1.51%  "
$ws.Range("D16").Value = "'0.00001108"
$ws.Range("E16").Value = "'  -0.95%  "
$ws.Range("D17").Value = "'91.20"
$ws.Range("E17").Value = "'  -0.77%  "
$ws.Range("D18").Value = "'0.06737"
$ws.Range("E18").Value = "'  -1.27%  "
$ws.Range("D19").Value = "'17.70"
$ws.Range("E19").Value = "'  -1.90%  "
$ws.Range("D20").Value = "'1.012"
$ws.Range("E20").Value = "'  -1.76%  "
$ws.Range("D21").Value = "'5.949"
$ws.Range("E21").Value = "'  -2.57%  "
$ws.Range("D22").Value = "'28.529.00"
$ws.Range("E22").Value = "'  -0.92%  "
$ws.Range("D23").Value = "'11.14"
$ws.Range("E23").Value = "'  -1.35%  "
$ws.Range("D24").Value = "'2.270"
$ws.Range("E24").Value = "'  -1.29%  "
$ws.Range("D25").Value = "'2.103.17"
$ws.Range("E25").Value = "'  -1.66%  "
$ws.Range("D26").Value = "'161.17"
$ws.Range("E26").Value = "'  -1.14%  "
$ws.Range("D27").Value = "'20.73"
$ws.Range("E27").Value = "'  -1.60%  "
$ws.Range("D28").Value = "'2.384"
$ws.Range("E28").Value = "'  -3.33%  "
$ws.Range("D29").Value = "'126.43"
$ws.Range("E29").Value = "'  -1.22%  "
$ws.Range("D30").Value = "'0.1055"
$ws.Range("E30").Value = "'  -0.76%  "
$ws.Range("D31").Value = "'1.048"
$ws.Range("E31").Value = "'  -0.92%  "
$ws.Range("D32").Value = "'5.786"
$ws.Range("E32").Value = "'  -3.36%  "
$ws.Range("D33").Value = "'3.610"
$ws.Range("E33").Value = "'  -2.13%  "
$ws.Range("D34").Value = "'0.02439"
$ws.Range("E34").Value = "'  -1.42%  "
$ws.Range("D35").Value = "'0.06502"
$ws.Range("E35").Value = "'  -2.55%  "
$ws.Range("D36").Value = "'0.2183"
$ws.Range("E36").Value = "'  -1.98%  "
$ws.Range("D37").Value = "'8.915"
$ws.Range("E37").Value = "'  -5.96%  "
$ws.Range("E38").Value = "'  +0.53%  "
$ws.Range("B39").Value = "'ARBITRUM"
$ws.Range("C39").Value = "'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D39").Value = "'1.190"
$ws.Range("E39").Value = "'  -1.02%  "
$ws.Range("B40").Value = "'TheSandbox"
$ws.Range("C40").Value = "'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D40").Value = "'0.6455"
$ws.Range("E40").Value = "'  -2.09%  "
$ws.Range("D41").Value = "'5.057"
$ws.Range("E41").Value = "'  +0.41%  "
$ws.Range("D42").Value = "'11.20"
$ws.Range("E42").Value = "'  +0.06%  "
$ws.Range("D43").Value = "'1.012"
$ws.Range("E43").Value = "'  -1.73%  "
$ws.Range("D44").Value = "'0.6071"
$ws.Range("E44").Value = "'  -1.89%  "
$ws.Range("E45").Value = "'  -1.02%  "
$ws.Range("D46").Value = "'3.709"
$ws.Range("E46").Value = "'  -1.46%  "
$ws.Range("D47").Value = "'2.020"
$ws.Range("E47").Value = "'  -0.33%  "
$ws.Range("D48").Value = "'1.203"
$ws.Range("E48").Value = "'  -7.73%  "
$ws.Range("D49").Value = "'1.212"
$ws.Range("E49").Value = "'  -2.58%  "
$ws.Range("D50").Value = "'122.32"
$ws.Range("E50").Value = "'  -0.66%  "
$ws.Range("D51").Value = "'0.06841"
$ws.Range("E51").Value = "'  -1.79%  "
